$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Characters(21, 2).Text = "51"
$ws.Range("C9").Characters(27, 10).Text = "12/19/2022"
$ws.Range("C9").Characters(48, 10).Text = "12/25/2022"

# --- Crime statistics table updates (rows 14-29) ---
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "0"
$ws.Range("N14").Value = -67.441860465116
$ws.Range("D15").Value = 1
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 32
$ws.Range("K15").Value = 15.625
$ws.Range("L15").Value = 48
$ws.Range("M15").Value = 37.037037037037
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 116.666666666667
$ws.Range("F16").Value = 34
$ws.Range("G16").Value = 36
$ws.Range("H16").Value = -5.555555555555
$ws.Range("I16").Value = 398
$ws.Range("J16").Value = 365
$ws.Range("K16").Value = 9.04109589041
$ws.Range("L16").Value = 16.715542521994
$ws.Range("M16").Value = -19.595959595959
$ws.Range("N16").Value = -76.670574443141
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -31.25
$ws.Range("G17").Value = 49
$ws.Range("H17").Value = -18.367346938775
$ws.Range("I17").Value = 591
$ws.Range("J17").Value = 545
$ws.Range("K17").Value = 8.440366972477
$ws.Range("L17").Value = 23.640167364016
$ws.Range("M17").Value = 33.710407239819
$ws.Range("N17").Value = -4.213938411669
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -60
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 7.692307692307
$ws.Range("I18").Value = 191
$ws.Range("J18").Value = 175
$ws.Range("K18").Value = 9.142857142857
$ws.Range("L18").Value = -17.672413793103
$ws.Range("M18").Value = -55.3738317757
$ws.Range("N18").Value = -91.763691246226
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 29
$ws.Range("E19").Value = -55.172413793103
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 96
$ws.Range("H19").Value = -54.166666666666
$ws.Range("I19").Value = 649
$ws.Range("J19").Value = 699
$ws.Range("K19").Value = -7.153075822603
$ws.Range("L19").Value = 8.892617449664
$ws.Range("M19").Value = 15.892857142857
$ws.Range("N19").Value = -28.993435448577
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 400
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 63.636363636363
$ws.Range("I20").Value = 249
$ws.Range("J20").Value = 210
$ws.Range("K20").Value = 18.571428571428
$ws.Range("L20").Value = 35.326086956521
$ws.Range("M20").Value = 35.326086956521
$ws.Range("N20").Value = -82.827586206896
$ws.Range("C21").Value = 44
$ws.Range("D21").Value = 58
$ws.Range("E21").Value = -24.137931034482
$ws.Range("F21").Value = 154
$ws.Range("G21").Value = 210
$ws.Range("H21").Value = -26.666666666666
$ws.Range("I21").Value = 2129
$ws.Range("J21").Value = 2033
$ws.Range("K21").Value = 4.722085587801
$ws.Range("L21").Value = 14.094319399785
$ws.Range("M21").Value = -0.838379133674
$ws.Range("N21").Value = -70.052046701364
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -50
$ws.Range("L22").Value = -11.538461538461
$ws.Range("M22").Value = 58.620689655172
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -26.666666666666
$ws.Range("F24").Value = 154
$ws.Range("G24").Value = 135
$ws.Range("H24").Value = 14.074074074074
$ws.Range("I24").Value = 2563
$ws.Range("J24").Value = 1531
$ws.Range("K24").Value = 67.406923579359
$ws.Range("L24").Value = 73.644986449864
$ws.Range("M24").Value = 103.25138778747
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -7.142857142857
$ws.Range("F25").Value = 70
$ws.Range("G25").Value = 69
$ws.Range("H25").Value = 1.449275362318
$ws.Range("I25").Value = 791
$ws.Range("J25").Value = 778
$ws.Range("K25").Value = 1.670951156812
$ws.Range("L25").Value = 3.942181340341
$ws.Range("M25").Value = -10.215664018161
$ws.Range("D26").Value = 3
$ws.Range("G26").Value = 10
$ws.Range("H26").Value = -40
$ws.Range("I26").Value = 72
$ws.Range("J26").Value = 50
$ws.Range("K26").Value = 44
$ws.Range("L26").Value = 80
$ws.Range("C27").Value = 2
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 80
$ws.Range("I27").Value = 90
$ws.Range("K27").Value = -15.094339622641
$ws.Range("L27").Value = 21.621621621621
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Value = 2
$ws.Range("E28").NumberFormat = "General"
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -60
$ws.Range("J28").Value = 42
$ws.Range("K28").Value = -14.285714285714
$ws.Range("N28").Value = -58.620689655172
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Value = 1
$ws.Range("E29").NumberFormat = "General"
$ws.Range("E29").Value = -100
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("J29").Value = 34
$ws.Range("K29").Value = -8.823529411764
$ws.Range("N29").Value = -61.728395061728
